# Kennard-Stone splitting: add a new "inv-12r-adj" column right before the
# existing "12r-adj" column (current column C), pushing every later column
# one slot to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new column at C; Excel shifts C:BR -> D:BS and carries the
# header's cell style (s="1") along with it, matching a native
# Insert-Column operation.
$ws.Columns("C").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "inv-12r-adj"

# New data values for the inserted column (rows 2-7).
$ws.Cells.Item(2, 3).Value = 0.36613973103532965
$ws.Cells.Item(3, 3).Value = 1.029044147065822
$ws.Cells.Item(4, 3).Value = 1.1688787796800191
$ws.Cells.Item(5, 3).Value = 0.22527225202136897
$ws.Cells.Item(6, 3).Value = 0.26774416639352816
$ws.Cells.Item(7, 3).Value = 0.24881160525941676

# Match the author's final selection: the newly added column's data range.
$ws.Range("C2:C7").Select()
